$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Recreate Q5RegSummary (delete + re-add => fresh sheetId, matching the
#    "sheetId 5 -> 7" change seen in workbook.xml; a brand-new sheet also has
#    no <cols> overrides, matching the target column-width reset).
#    NOTE: populate all data on the new sheet BEFORE moving it back into
#    position - moving first causes the object handle to rebind to whatever
#    sheet ends up at the old numeric slot.
# ---------------------------------------------------------------------------
$oldQ5 = $wb.Worksheets.Item("Q5RegSummary")
$oldQ5.Delete()

$q5 = $wb.Worksheets.Add()
$q5.Name = "Q5RegSummary"

$q5.Range("B1").Value = "term"
$q5.Range("C1").Value = "estimate"
$q5.Range("D1").Value = "std.error"
$q5.Range("E1").Value = "statistic"
$q5.Range("F1").Value = "p.value"

$q5.Range("A2").Value = 1
$q5.Range("B2").Value = "(Intercept)"
$q5.Range("C2").Value = -1.9702928229901899
$q5.Range("D2").Value = 0.54531605123090998
$q5.Range("E2").Value = -3.6131209021681299
$q5.Range("F2").Value = 0.00036385163339560002

$q5.Range("A3").Value = 2
$q5.Range("B3").Value = "LnPrice"
$q5.Range("C3").Value = -0.0057904941437127997
$q5.Range("D3").Value = 0.027748336049512701
$q5.Range("E3").Value = -0.20867896847510201
$q5.Range("F3").Value = 0.83486415545924098

$q5.Range("A4").Value = 3
$q5.Range("B4").Value = "LnPrint"
$q5.Range("C4").Value = 0.018269009269862801
$q5.Range("D4").Value = 0.0046180295795582096
$q5.Range("E4").Value = 3.9560182443895302
$q5.Range("F4").Value = 0.000098633033965413004
$q5.Range("F4").NumberFormat = "0.00E+00"

$q5.Range("A5").Value = 4
$q5.Range("B5").Value = "LnOut"
$q5.Range("C5").Value = -0.0067009009023669399
$q5.Range("D5").Value = 0.0059209043840908002
$q5.Range("E5").Value = -1.13173604363076
$q5.Range("F5").Value = 0.25880030987978397

$q5.Range("A6").Value = 5
$q5.Range("B6").Value = "LnBroad"
$q5.Range("C6").Value = 0.0035672735736858598
$q5.Range("D6").Value = 0.0057301758659404102
$q5.Range("E6").Value = 0.62254172596854696
$q5.Range("F6").Value = 0.53413768327368305

$q5.Range("A7").Value = 6
$q5.Range("B7").Value = "LagTotalMinusSales"
$q5.Range("C7").Value = 0.000031969885502778801
$q5.Range("D7").Value = 0.0000092606167482020302
$q5.Range("E7").Value = 3.4522415052956199
$q5.Range("F7").Value = 0.00064989707420020498
$q5.Range("C7:D7").NumberFormat = "0.00E+00"

$q5.Columns.Item(1).ColumnWidth = 2.83

# put Q5RegSummary back between Q4RegSummary and Q6RegSummary
$q6 = $wb.Worksheets.Item("Q6RegSummary")
$q5.Move($q6)

# ---------------------------------------------------------------------------
# 2) Recreate Q3RegSummary the same way (sheetId 4 -> 8). Doing this *after*
#    Q5RegSummary reproduces the sheetId allocation order seen in the diff
#    (Q5 -> 7, Q3 -> 8).
# ---------------------------------------------------------------------------
$oldQ3 = $wb.Worksheets.Item("Q3RegSummary")
$oldQ3.Delete()

$q3 = $wb.Worksheets.Add()
$q3.Name = "Q3RegSummary"

$q3.Range("B1").Value = "term"
$q3.Range("C1").Value = "estimate"
$q3.Range("D1").Value = "std.error"
$q3.Range("E1").Value = "statistic"
$q3.Range("F1").Value = "p.value"

$q3.Range("A2").Value = 1
$q3.Range("B2").Value = "(Intercept)"
$q3.Range("C2").Value = -0.11736535028180101
$q3.Range("D2").Value = 0.098395868188021304
$q3.Range("E2").Value = -1.1927873847053401
$q3.Range("F2").Value = 0.23404845575907099

$q3.Range("A3").Value = 2
$q3.Range("B3").Value = "LnPrice"
$q3.Range("C3").Value = 0.036774462363971298
$q3.Range("D3").Value = 0.025379123556186999
$q3.Range("E3").Value = 1.4490044261204
$q3.Range("F3").Value = 0.14855067433341401

$q3.Range("A4").Value = 3
$q3.Range("B4").Value = "LnPrint"
$q3.Range("C4").Value = 0.014776084997315199
$q3.Range("D4").Value = 0.0046001787421577796
$q3.Range("E4").Value = 3.21206757944024
$q3.Range("F4").Value = 0.0014851890583199999

$q3.Range("A5").Value = 4
$q3.Range("B5").Value = "LnOut"
$q3.Range("C5").Value = -0.01262205486868
$q3.Range("D5").Value = 0.0057857043577503302
$q3.Range("E5").Value = -2.1815934738822098
$q3.Range("F5").Value = 0.0300423405364139

$q3.Range("A6").Value = 5
$q3.Range("B6").Value = "LnBroad"
$q3.Range("C6").Value = -0.0057637188827885204
$q3.Range("D6").Value = 0.0051584593239443303
$q3.Range("E6").Value = -1.1173333975969799
$q3.Range("F6").Value = 0.26489119983752701

$q3.Columns.Item(1).ColumnWidth = 2.83

# put Q3RegSummary back between Q2RegSummary and Q4RegSummary
$q4 = $wb.Worksheets.Item("Q4RegSummary")
$q3.Move($q4)

# ---------------------------------------------------------------------------
# 3) Active-tab move: Q6RegSummary (index 4) was the active / tabSelected
#    sheet; the target has no explicit active tab override and Q2RegSummary
#    picks up tabSelected instead, so activate the first sheet.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("Q2RegSummary")
$q2.Activate()
